# Dudas y anotaciones.docx - añadir contenido al final de la memoria del TFM
$d = $word.ActiveDocument

# --- Paragraph 1: "Preguntar si es necesario seguir un orden cronológico..." ---
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Preguntar si es necesario seguir un orden cronológico en los artículos mencionados."

# --- Paragraph 2: "OLIN -> ..." ---
$r = $p1.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "OLIN -> The connectionist nature of the info-fuzzy network (each terminal node is connected to every target node) resembles the topological structure of multi-layer neural networks (see [27]), which also have input and output nodes and a variable number of hidden layers. Consequently, we define our model as a network and not as a tree. -> No obstante, lo compara con el CVFDT -> ¿Lo pongo en las propuestas de árboles de decisión?"

# --- Paragraph 3: "UFFT -> ..." with the _GoBack bookmark relocated mid-paragraph ---
$r = $p2.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "UFFT -> Para un problema multiclase construye un bosque de árboles de decisión binario, uno para cada par de valores que puede tomar la variable clase"

# Move the existing _GoBack bookmark from the previous paragraph to right after
# "...variable clase" in this new paragraph, then append the remaining text.
$gb = $d.Bookmarks("\_GoBack")
$gb.Delete()

$r3 = $p3.Range
$r3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r3) | Out-Null

$r3b = $d.Range($p3.Range.End, $p3.Range.End)
$r3b.InsertAfter(". ¿Es cierto que no es ensemble?")

# --- Paragraph 4: trailing empty paragraph ---
$r = $p3.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.ListFormat.RemoveNumbers()
$p4.Style = $d.Styles("Normal")
